$wb = $excel.ActiveWorkbook

# Remove the "Hướng dẫn" (guide) worksheet entirely.
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Hướng dẫn").Delete()

# Work on the remaining "Danh sách" sheet.
$ws = $wb.Worksheets.Item("Danh sách")

# Update the employee code on row 2 from "00517" to "L0001".
$ws.Range("A2").Value = "L0001"

# Remove row 3 (the second sample data row), leaving only the header + one row.
$ws.Rows.Item(3).Delete()
